# Updated cryptos list refresh (prices / 1h volume deltas), plus two
# rank swaps (VeChain <-> TrustWalletToken, EnergySwap <-> PaxDollar).
#
# Price cells in column D are plain text in the source sheet (e.g.
# "311.56", "1.832.22"). Values that parse as a single plain number
# (one decimal point, no thousands separators) would otherwise be
# auto-converted to a numeric cell by the normal Value-assignment
# type inference, so those are entered with a leading apostrophe
# (text/quote prefix) the same way a user would force text entry,
# then the cell style is reset back to Normal so no stray number
# format lingers on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($Address, $Value) {
    # Plain decimal numbers (e.g. "311.70", "0.07830") would otherwise be
    # silently reinterpreted as numeric values (dropping trailing zeros /
    # text-vs-number formatting) by ordinary Value assignment. Values with
    # more than one '.' (e.g. "27.083.95") or any non-numeric character
    # are already safe as-is.
    $needsPrefix = ($Value -match '^[+-]?\d+(\.\d+)?$')
    $range = $ws.Range($Address)
    if ($needsPrefix) {
        $range.Value = "'" + $Value
        $range.Style = "Normal"
    } else {
        $range.Value = $Value
    }
}

# row -> column -> new value, taken from the diff
$rowEdits = [ordered]@{
    2  = @{ D = "27.083.95";    E = "  -0.04%  " }
    3  = @{ D = "1.833.79";     E = "  +0.46%  " }
    4  = @{ D = "1.008";        E = "  -0.01%  " }
    5  = @{ D = "311.70" }
    6  = @{ D = "1.006";        E = "  -0.09%  " }
    7  = @{ D = "0.4651";       E = "  -0.71%  " }
    8  = @{ D = "0.3710";       E = "  +1.68%  " }
    9  = @{ D = "0.07367";      E = "  -0.24%  " }
    10 = @{ D = "0.8751";       E = "  -0.39%  " }
    11 = @{ D = "20.01";        E = "  -1.20%  " }
    12 = @{ D = "0.07830";      E = "  +4.59%  " }
    13 = @{ E = "  +1.61%  " }
    14 = @{ D = "5.353";        E = "  -0.30%  " }
    15 = @{ D = "91.98";        E = "  -1.02%  " }
    16 = @{ D = "1.712.92";     E = "  -9.63%  " }
    17 = @{ E = "  +0.25%  " }
    18 = @{ D = "0.000008856";  E = "  +1.60%  " }
    19 = @{ D = "1.008";        E = "  +0.10%  " }
    20 = @{ D = "27.440.32" }
    21 = @{ E = "  +0.31%  " }
    22 = @{ D = "5.143";        E = "  -1.70%  " }
    23 = @{ E = "  -0.31%  " }
    24 = @{ D = "1.931.91";     E = "  -7.12%  " }
    25 = @{ D = "152.31";       E = "  +0.60%  " }
    26 = @{ D = "1.826";        E = "  -2.98%  " }
    27 = @{ D = "18.34";        E = "  -0.85%  " }
    28 = @{ D = "2.100" }
    29 = @{ D = "5.085";        E = "  -1.48%  " }
    30 = @{ D = "115.49";       E = "  -0.72%  " }
    31 = @{ D = "0.08878";      E = "  -0.48%  " }
    32 = @{ D = "2.962";        E = "  +0.68%  " }
    33 = @{ D = "0.7292";       E = "  -1.94%  " }
    34 = @{ D = "4.446";        E = "  -1.35%  " }
    35 = @{ D = "1.140";        E = "  -1.82%  " }
    36 = @{ D = "2.482";        E = "  -2.06%  " }
    37 = @{ B = "VeChain";           C = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet";          D = "0.01950"; E = "  +0.97%  " }
    38 = @{ B = "TrustWalletToken";  C = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt";     D = "1.070";   E = "  -1.84%  " }
    39 = @{ D = "0.05222";      E = "  -1.27%  " }
    40 = @{ D = "2.931";        E = "  -0.11%  " }
    41 = @{ D = "7.179" }
    42 = @{ D = "0.5202";       E = "  -0.83%  " }
    43 = @{ D = "0.8835";       E = "  -12.26%  " }
    44 = @{ D = "0.1631";       E = "  -0.53%  " }
    45 = @{ D = "8.234";        E = "  -1.58%  " }
    46 = @{ D = "0.4830";       E = "  -1.30%  " }
    47 = @{ B = "EnergySwap";  C = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; D = "10.23"; E = "  -1.66%  " }
    48 = @{ B = "PaxDollar";   C = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"; D = "1.006"; E = "  -0.13%  " }
    49 = @{ D = "102.80";      E = "  -1.48%  " }
    50 = @{ D = "1.629";       E = "  -1.19%  " }
    51 = @{ D = "0.06224";     E = "  -0.70%  " }
}

foreach ($row in $rowEdits.Keys) {
    $cols = $rowEdits[$row]
    foreach ($col in @("B", "C", "D", "E")) {
        if ($cols.Contains($col)) {
            $address = "$col$row"
            $newValue = $cols[$col]
            Set-TextCell $address $newValue
        }
    }
}
